$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("A16").Value = '2025-04-18 10:25'
$ws.Range("B16").Value = 'http://www.bzqzf.gov.cn/group2/M00/06/EB/rBUtIWSvnuGAAzraAAK6ABoS7RQ356.xls'
$ws.Range("C16").Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=4232bdec735d418898416f8f652f2152&type=0'
$ws.Range("D16").Value = '“三公经费”'
$ws.Range("E16").Value = '“三公”经费'
$ws.Range("F16").Value = '巴州区老干部局2016年预算批复表.xls'
$ws.Range("G16").Value = 'http://www.bzqzf.gov.cn/zwgk/zdxxgk/czxx/czyjs/8433911.html'
$ws.Range("H16").Value = '中共巴中市巴州区委老干部局2016年部门预算编制说明'

# Row 17
$ws.Range("A17").Value = '2025-04-18 10:25'
$ws.Range("B17").Value = 'http://www.bzqzf.gov.cn/oldfiles/bzq/upload/59a767e47f8b9acf685d2650/201613/201613_57099b7504b32a646c786a8ef778a92e.xls?fileName=57099b7504b32a646c786a8ef778a92e.xls'
$ws.Range("C17").Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=a555ceb4faf44a4296d577fbacff3d05&type=0'
$ws.Range("D17").Value = '“三公经费”'
$ws.Range("E17").Value = '“三公”经费'
$ws.Range("F17").Value = '附件：57099b7504b32a646c786a8ef778a92e.xls'
$ws.Range("G17").Value = 'http://www.bzqzf.gov.cn/zwgk/zdxxgk/czxx/czyjs/8432891.html'
$ws.Range("H17").Value = '巴中市巴州区住建系统2017年部门预算编制说明'

# Row 18
$ws.Range("A18").Value = '2025-04-18 10:25'
$ws.Range("B18").Value = 'http://www.bzqzf.gov.cn/oldfiles/bzq/upload/59a767e47f8b9acf685d2650/201613/201613_3efba75a0c0f7184f8826379e0ceb648.xls?fileName=3efba75a0c0f7184f8826379e0ceb648.xls'
$ws.Range("C18").Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=fd8ff0b4a5ce4e9287ba0ef57607275f&type=0'
$ws.Range("D18").Value = '“三公经费”'
$ws.Range("E18").Value = '“三公”经费'
$ws.Range("F18").Value = '附件：3efba75a0c0f7184f8826379e0ceb648.xls'
$ws.Range("G18").Value = 'http://www.bzqzf.gov.cn/zwgk/zdxxgk/czxx/czyjs/8433621.html'
$ws.Range("H18").Value = '巴中市巴州区青少年宫2017年部门预算编制说明'

# Row 19
$ws.Range("A19").Value = '2025-04-18 10:25'
$ws.Range("B19").Value = 'http://www.bzqzf.gov.cn/oldfiles/bzq/upload/59a767e47f8b9acf685d2650/201613/201613_436db511abd2661d742920515de4b318.xls?fileName=436db511abd2661d742920515de4b318.xls'
$ws.Range("C19").Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=dcadb047be444a0c80e3009a2f437f6d&type=0'
$ws.Range("D19").Value = '“三公经费”'
$ws.Range("E19").Value = '“三公”经费'
$ws.Range("F19").Value = '附件：436db511abd2661d742920515de4b318.xls'
$ws.Range("G19").Value = 'http://www.bzqzf.gov.cn/zwgk/zdxxgk/czxx/czyjs/8432911.html'
$ws.Range("H19").Value = '巴中市巴州区经济和信息化局2017年部门预算编制说明'

